$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B26: clarify that all fields are required when adding a product ---
$ws.Range("B26").Value = "Supplier can add product with data (product id,product photo, product price, product version, product platform [IOS&Android] ) all fileds are required"

# --- New supplier requirement rows 32-34 (previously blank placeholder rows) ---
# Copy the formatting from the existing supplier-table rows (29-31) so the
# new rows 32-34 match the table's look (borders/fonts) for columns A and C,
# while column B keeps its pre-existing format (only its value changes).
$ws.Range("A31").Copy()
$ws.Range("A32").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("C31").Copy()
$ws.Range("C32").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A31").Copy()
$ws.Range("A33").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("C31").Copy()
$ws.Range("C33").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A31").Copy()
$ws.Range("A34").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("C31").Copy()
$ws.Range("C34").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = $false

$ws.Range("A32").Value = "CRS_supplier_007"
$ws.Range("B32").Value = "The error message, in red color,  appears in the top of the form of add product and update product "
$ws.Range("C32").Value = "CRS-SIQ_supplier_007"

$ws.Range("A33").Value = "CRS_supplier_008"
$ws.Range("B33").Value = "The first filed that has an error will be highlighted by red color and exclamation mark appears in the right of the field "
$ws.Range("C33").Value = "CRS-SIQ_supplier_008"

$ws.Range("A34").Value = "CRS_supplier_009"
$ws.Range("B34").Value = "The product id should be unique, there are not two products with the same product id"
$ws.Range("C34").Value = "CRS-SIQ_supplier_009"

# Match the row heights used by the rest of the supplier table (18pt), and
# keep row 34's pre-existing fixed height (15.75pt).
$ws.Rows.Item(32).RowHeight = 18
$ws.Rows.Item(33).RowHeight = 18
$ws.Rows.Item(34).RowHeight = 15.75

# --- View state: zoom out a bit and move the selection down to the new rows ---
$excel.ActiveWindow.Zoom = 84
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B36").Select()
